# Update the "想去人数" (interested-people count) figures that were refreshed
# in this scrape run. Two worksheets share the same underlying data set:
#   sheet1 ("展览")   -> Worksheets.Item(1)
#   sheet4 ("全部类型") -> Worksheets.Item(4)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F7").Value  = 56
$ws1.Range("F13").Value = 774
$ws1.Range("F14").Value = 6
$ws1.Range("F15").Value = 8
$ws1.Range("F16").Value = 1534
$ws1.Range("F17").Value = 1534
$ws1.Range("F18").Value = 898
$ws1.Range("F19").Value = 32
$ws1.Range("F22").Value = 352
$ws1.Range("F26").Value = 6723
$ws1.Range("F27").Value = 5099
$ws1.Range("F28").Value = 5099
$ws1.Range("F31").Value = 211
$ws1.Range("F41").Value = 266

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F11").Value = 56
$ws4.Range("F17").Value = 774
$ws4.Range("F18").Value = 1534
$ws4.Range("F19").Value = 1534
$ws4.Range("F20").Value = 898
$ws4.Range("F21").Value = 32
$ws4.Range("F24").Value = 352
$ws4.Range("F29").Value = 6723
$ws4.Range("F30").Value = 5099
$ws4.Range("F31").Value = 5099
$ws4.Range("F45").Value = 266

$wb.Save()
